# Insert a new row for "mistral_7b_instruct_v2" right before the existing
# "mistral_7b_instruct_v3" row (row 8), pushing all following rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("8:8").Insert()

$ws.Range("A8").Value = "mistral_7b_instruct_v2"
$ws.Range("B8").Value = 2486
$ws.Range("C8").Value = 1059
$ws.Range("D8").Value = 1295
$ws.Range("E8").Value = 1070
$ws.Range("F8").Value = 155
$ws.Range("G8").Value = 107
$ws.Range("H8").Value = 66
$ws.Range("I8").Value = 72
